$d = $word.ActiveDocument

# The document currently reads "...targeting iOS6.0 or higher." and has a
# stray _GoBack bookmark sitting by itself in the empty paragraph just
# above that sentence. The edit bumps the minimum supported iOS version
# down to 5.x, and (per the target markup) the _GoBack bookmark ends up
# relocated to sit right after "iOS5" / before ".0 or higher." in that
# same sentence.

# 1) Drop the old _GoBack bookmark from wherever it currently lives.
If ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) "iOS6.0 or higher." -> "iOS5.0 or higher."
$rFind = $d.Content
$rFind.Find.Execute("iOS6", $true, $false, $false, $false, $false, $true, 1, $false, "iOS5", 2) | Out-Null

# 3) Re-create the _GoBack bookmark as a collapsed range right after the
#    newly written "iOS5", splitting the run so ".0 or higher." starts a
#    new run just past the bookmark.
$rMark = $d.Content
$rMark.Find.Execute("iOS5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rMark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rMark) | Out-Null
